$d = $word.ActiveDocument

# "3.37" becomes "3.27" in the problem list. Locate the run with a Find so we
# know exactly where the replacement text lands, then overwrite its text.
$match = $d.Content
$found = $match.Find.Execute("3.37", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$matchStart = $match.Start
$match.Text = "3.27"

# The document's auto "_GoBack" bookmark currently sits right after "4.11"
# (Word drops it at the last edited spot). In the revised document it has
# moved into the middle of the new "3.27" run, splitting it into "3.2" + the
# (empty) bookmark + "7". Recreate that: drop the old bookmark and add a new
# collapsed one 3 characters into the replacement text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$splitPoint = $matchStart + 3
$bookmarkRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
